$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 6).Value = 5841
$ws.Cells.Item(2, 10).Value = 889
$ws.Cells.Item(3, 9).Value = 7487
$ws.Cells.Item(3, 10).Value = 978
$ws.Cells.Item(4, 8).Value = 1685
$ws.Cells.Item(4, 9).Value = 1752
$ws.Cells.Item(4, 10).Value = 222
$ws.Cells.Item(5, 10).Value = 72
$ws.Cells.Item(6, 10).Value = 1370
$ws.Cells.Item(7, 6).Value = 24070
$ws.Cells.Item(7, 8).Value = 25997
$ws.Cells.Item(7, 10).Value = 3531

$ws = $wb.Worksheets.Item("Uptown")
$ws.Cells.Item(2, 10).Value = 13
$ws.Cells.Item(6, 10).Value = 14
$ws.Cells.Item(7, 10).Value = 42

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Cells.Item(3, 10).Value = 12
$ws.Cells.Item(4, 10).Value = 5
$ws.Cells.Item(7, 10).Value = 44

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(2, 10).Value = 33
$ws.Cells.Item(3, 10).Value = 40
$ws.Cells.Item(7, 10).Value = 121

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Cells.Item(6, 10).Value = 15
$ws.Cells.Item(7, 10).Value = 47

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Cells.Item(2, 10).Value = 14
$ws.Cells.Item(7, 10).Value = 25

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(5, 10).Value = 11
$ws.Cells.Item(8, 10).Value = 226
$ws.Cells.Item(10, 10).Value = 23
$ws.Cells.Item(19, 10).Value = 114
$ws.Cells.Item(20, 10).Value = 80
$ws.Cells.Item(21, 10).Value = 8
$ws.Cells.Item(23, 10).Value = 32
$ws.Cells.Item(27, 10).Value = 17
$ws.Cells.Item(29, 9).Value = 1554
$ws.Cells.Item(29, 10).Value = 182
$ws.Cells.Item(31, 10).Value = 25
$ws.Cells.Item(32, 10).Value = 8
$ws.Cells.Item(33, 10).Value = 150
$ws.Cells.Item(36, 10).Value = 53
$ws.Cells.Item(37, 10).Value = 121
$ws.Cells.Item(47, 10).Value = 32
$ws.Cells.Item(50, 10).Value = 19
$ws.Cells.Item(51, 10).Value = 46
$ws.Cells.Item(52, 10).Value = 78
$ws.Cells.Item(54, 10).Value = 63
$ws.Cells.Item(63, 6).Value = 171
$ws.Cells.Item(63, 9).Value = 176
$ws.Cells.Item(71, 10).Value = 16
$ws.Cells.Item(72, 10).Value = 14
$ws.Cells.Item(76, 10).Value = 57
$ws.Cells.Item(77, 8).Value = 161
$ws.Cells.Item(77, 10).Value = 32
$ws.Cells.Item(78, 10).Value = 48
$ws.Cells.Item(81, 10).Value = 4
$ws.Cells.Item(85, 10).Value = 144
$ws.Cells.Item(89, 10).Value = 42
$ws.Cells.Item(90, 10).Value = 42
$ws.Cells.Item(91, 10).Value = 45
$ws.Cells.Item(94, 10).Value = 20
$ws.Cells.Item(95, 10).Value = 66
$ws.Cells.Item(96, 10).Value = 44
$ws.Cells.Item(99, 10).Value = 47
$ws.Cells.Item(101, 6).Value = 24070
$ws.Cells.Item(101, 8).Value = 25997
$ws.Cells.Item(101, 10).Value = 3531

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Cells.Item(2, 10).Value = 24
$ws.Cells.Item(4, 10).Value = 2
$ws.Cells.Item(7, 10).Value = 66

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(6, 10).Value = 67
$ws.Cells.Item(7, 10).Value = 150

$ws = $wb.Worksheets.Item("Loop")
$ws.Cells.Item(2, 10).Value = 19
$ws.Cells.Item(7, 10).Value = 63

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(2, 10).Value = 54
$ws.Cells.Item(3, 9).Value = 530
$ws.Cells.Item(4, 10).Value = 7
$ws.Cells.Item(7, 9).Value = 1554
$ws.Cells.Item(7, 10).Value = 182

$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(2, 10).Value = 22
$ws.Cells.Item(7, 10).Value = 114

$ws = $wb.Worksheets.Item("River North")
$ws.Cells.Item(2, 10).Value = 6
$ws.Cells.Item(3, 10).Value = 12
$ws.Cells.Item(6, 10).Value = 35
$ws.Cells.Item(7, 10).Value = 57

$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(3, 10).Value = 47
$ws.Cells.Item(7, 10).Value = 144

$ws = $wb.Worksheets.Item("Avondale")
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(6, 10).Value = 9
$ws.Cells.Item(7, 10).Value = 23

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Cells.Item(3, 10).Value = 18
$ws.Cells.Item(4, 10).Value = 6
$ws.Cells.Item(7, 10).Value = 48

$ws = $wb.Worksheets.Item("Douglas")
$ws.Cells.Item(3, 10).Value = 8
$ws.Cells.Item(7, 10).Value = 32

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Cells.Item(3, 10).Value = 19
$ws.Cells.Item(7, 10).Value = 45

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(7, 10).Value = 8

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(2, 10).Value = 20
$ws.Cells.Item(3, 10).Value = 27
$ws.Cells.Item(7, 10).Value = 80

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Cells.Item(2, 10).Value = 15
$ws.Cells.Item(7, 10).Value = 53

$ws = $wb.Worksheets.Item("Little Village")
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(7, 10).Value = 78

$ws = $wb.Worksheets.Item("West Loop")
$ws.Cells.Item(6, 10).Value = 11
$ws.Cells.Item(7, 10).Value = 20

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Cells.Item(2, 10).Value = 6
$ws.Cells.Item(7, 10).Value = 32

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Cells.Item(2, 10).Value = 4
$ws.Cells.Item(4, 10).Value = 4
$ws.Cells.Item(7, 10).Value = 19

$ws = $wb.Worksheets.Item("Galewood")
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(7, 10).Value = 8

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(2, 10).Value = 69
$ws.Cells.Item(7, 10).Value = 226

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Cells.Item(6, 10).Value = 7
$ws.Cells.Item(7, 10).Value = 11

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(7, 10).Value = 17

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Cells.Item(2, 10).Value = 10
$ws.Cells.Item(3, 10).Value = 11
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(7, 10).Value = 42

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Cells.Item(3, 10).Value = 17
$ws.Cells.Item(6, 10).Value = 7
$ws.Cells.Item(7, 10).Value = 46

$ws = $wb.Worksheets.Item("Oakland")
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(7, 10).Value = 16

$ws = $wb.Worksheets.Item("Old Town")
$ws.Cells.Item(2, 10).Value = 5
$ws.Cells.Item(7, 10).Value = 14

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Cells.Item(3, 10).Value = 9
$ws.Cells.Item(4, 8).Value = 6
$ws.Cells.Item(7, 8).Value = 161
$ws.Cells.Item(7, 10).Value = 32

$ws = $wb.Worksheets.Item("Sauganash,Forest Glen")
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(6, 10).Value = 4
